# Change #2 in word-doc:
# After the existing "Testing changes in word-document." paragraph, add an
# empty paragraph followed by a new paragraph containing
# "Hello! I made the second change in my test word document."
# The _GoBack bookmark (which marked the last edit position, right after
# the original sentence) ends up wrapping the end of the newly typed text,
# i.e. it moves from the end of paragraph 1 to the end of paragraph 3.

$d = $word.ActiveDocument

# The matching run formatting used throughout this document.
$rPrXml = "<w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/><w:sz w:val='28'/><w:szCs w:val='28'/><w:lang w:val='en-US'/></w:rPr>"

# Drop the old "_GoBack" bookmark sitting at the end of paragraph 1 -- it
# will be re-created at the end of the newly-typed text below.
$d.Bookmarks("_GoBack").Delete() | Out-Null

# Locate the end of the document's text (just before where the bookmark
# used to sit) and insert the two new paragraphs there, leaving paragraph 1
# untouched.
$endRange = $d.Content
$endRange.Collapse(0)

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$emptyParaXml = "<w:p $wNs><w:pPr>$rPrXml</w:pPr></w:p>"

$secondParaXml = "<w:p $wNs><w:pPr>$rPrXml</w:pPr><w:r>$rPrXml<w:t>Hello! I made the second change in my test word document.</w:t></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>"

$endRange.InsertXML($emptyParaXml + $secondParaXml) | Out-Null

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
